$d = $word.ActiveDocument

# --- Add an explanatory "first page" header with two right-aligned lines ---
$sec = $d.Sections(1)
$sec.PageSetup.DifferentFirstPageHeaderFooter = $true

$hdr = $sec.Headers(2)   # wdHeaderFooterFirstPage
$r = $hdr.Range

$r.InsertAfter("Line 1")
$r.Collapse(0)
$r.InsertParagraphAfter()
$hdr.Range.Paragraphs(2).Style = "Normal"
$r.Collapse(0)
$r.InsertAfter("Line 2")

$p1 = $hdr.Range.Paragraphs(1)
$p1.Alignment = 2   # wdAlignParagraphRight
$p2 = $hdr.Range.Paragraphs(2)
$p2.Alignment = 2   # wdAlignParagraphRight
$p2.SpaceAfter = 24

# --- Give the Normal style an explicit Arial font ---
$s = $d.Styles("Normal")
$s.Font.Name = "Arial"
$s.Font.NameFarEast = "Arial"
$s.Font.NameAscii = "Arial"
$s.Font.NameOther = "Arial"
$s.Font.NameBi = "Arial"
